$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" '42.633.06'
Set-TextCell "E2" '  -0.48%  '
Set-TextCell "D3" '2.295.53'
Set-TextCell "E3" '  -0.17%  '
Set-TextCell "E4" '  -0.03%  '
Set-TextCell "D5" '301.11'
Set-TextCell "E5" '  -1.59%  '
Set-TextCell "D6" '95.67'
Set-TextCell "E6" '  -1.19%  '
Set-TextCell "E7" '  -0.51%  '
Set-TextCell "E8" '  +0.05%  '
Set-TextCell "E9" '  -1.85%  '
Set-TextCell "D10" '34.56'
Set-TextCell "E11" '  +4.63%  '
Set-TextCell "E12" '  -0.93%  '
Set-TextCell "E13" '  -0.33%  '
Set-TextCell "E14" '  +0.08%  '
Set-TextCell "D15" '2.650.98'
Set-TextCell "E15" '  -0.31%  '
Set-TextCell "D16" '2.292.81'
Set-TextCell "E16" '  -0.19%  '
Set-TextCell "D17" '0.783'
Set-TextCell "D18" '42.550.46'
Set-TextCell "E18" '  -0.54%  '
Set-TextCell "D19" '12.30'
Set-TextCell "E19" '  -5.52%  '
Set-TextCell "E20" '  -1.09%  '
Set-TextCell "E21" '  -0.58%  '
Set-TextCell "D22" '67.76'
Set-TextCell "D23" '2.26'
Set-TextCell "E23" '  +5.39%  '
Set-TextCell "D24" '234.93'
Set-TextCell "E24" '  -0.52%  '
Set-TextCell "D26" '2.40'
Set-TextCell "E26" '  -2.77%  '
Set-TextCell "D27" '24.51'
Set-TextCell "E27" '  -3.41%  '
Set-TextCell "D28" '2.36'
Set-TextCell "E28" '  +14.70%  '
Set-TextCell "D29" '164.68'
Set-TextCell "E29" '  -0.89%  '
Set-TextCell "E30" '  -0.21%  '
Set-TextCell "D31" '32.14'
Set-TextCell "E31" '  -3.12%  '
Set-TextCell "E32" '  -0.06%  '
Set-TextCell "E33" '  -0.43%  '
Set-TextCell "E34" '  -1.54%  '
Set-TextCell "D35" '4.46'
Set-TextCell "E35" '  -6.64%  '
Set-TextCell "D36" '0.0702'
Set-TextCell "E36" '  +1.29%  '
Set-TextCell "E37" '  -3.01%  '
Set-TextCell "E38" '  -1.20%  '
Set-TextCell "E39" '  -0.51%  '
Set-TextCell "B40" 'Stellar'
Set-TextCell "C40" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D40" '0.108'
Set-TextCell "E40" '  -1.21%  '
Set-TextCell "B41" 'LidoDAOToken'
Set-TextCell "C41" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell "D41" '2.70'
Set-TextCell "E41" '  -0.35%  '
Set-TextCell "D42" '20.32'
Set-TextCell "E42" '  +11.25%  '
Set-TextCell "D43" '1.964.09'
Set-TextCell "E43" '  -1.96%  '
Set-TextCell "D44" '10.50'
Set-TextCell "E44" '  +5.11%  '
Set-TextCell "E46" '  -2.19%  '
Set-TextCell "E47" '  -0.41%  '
Set-TextCell "E48" '  +0.03%  '
Set-TextCell "D49" '2.520.64'
Set-TextCell "E49" '  -0.26%  '
Set-TextCell "D50" '53.17'
Set-TextCell "E50" '  -0.93%  '
Set-TextCell "D51" '71.27'
Set-TextCell "E51" '  -0.40%  '

Write-Output "Applied 77 cell updates"
